$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.346.14'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '3.018.41'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''594.89'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('D6').Value = '''149.08'
$ws.Range('E6').Value = '  +2.34%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.017.35'
$ws.Range('E8').Value = '  +0.63%  '
$ws.Range('E9').Value = '  -1.12%  '
$ws.Range('D10').Value = '''6.41'
$ws.Range('E10').Value = '  +10.99%  '
$ws.Range('E11').Value = '  +1.89%  '
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('E13').Value = '  +1.87%  '
$ws.Range('D14').Value = '''34.50'
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('D15').Value = '''0.126'
$ws.Range('E15').Value = '  +2.48%  '
$ws.Range('D16').Value = '3.517.63'
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('E17').Value = '  -0.54%  '
$ws.Range('D18').Value = '62.253.26'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('D19').Value = '3.018.67'
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('D20').Value = '''448.87'
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('E21').Value = '  +1.68%  '
$ws.Range('D22').Value = '''0.690'
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('D24').Value = '''82.33'
$ws.Range('E24').Value = '  +0.62%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').Value = '''2.25'
$ws.Range('E25').Value = '  +2.71%  '
$ws.Range('B26').Value = 'RenderToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D26').Value = '''10.88'
$ws.Range('E26').Value = '  +11.80%  '
$ws.Range('D27').Value = '''12.04'
$ws.Range('E27').Value = '  -1.34%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('E29').Value = '  +2.79%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').Value = '''7.19'
$ws.Range('E31').Value = '  +3.75%  '
$ws.Range('E32').Value = '  +2.28%  '
$ws.Range('D33').Value = '''27.54'
$ws.Range('E33').Value = '  -0.39%  '
$ws.Range('E34').Value = '  +1.13%  '
$ws.Range('D35').Value = '0.0₃0851'
$ws.Range('E35').Value = '  +6.16%  '
$ws.Range('E36').Value = '  +0.66%  '
$ws.Range('E37').Value = '  +2.00%  '
$ws.Range('D38').Value = '''2.07'
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('D39').Value = '''3.02'
$ws.Range('E39').Value = '  +4.99%  '
$ws.Range('D40').Value = '''50.13'
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('E41').Value = '  -1.21%  '
$ws.Range('E42').Value = '  +0.55%  '
$ws.Range('B43').Value = 'Arweave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D43').Value = '''41.45'
$ws.Range('E43').Value = '  +10.92%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').Value = '''0.285'
$ws.Range('E44').Value = '  +6.79%  '
$ws.Range('D45').Value = '''394.03'
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('E46').Value = '  -1.00%  '
$ws.Range('D47').Value = '2.739.79'
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('D48').Value = '''134.65'
$ws.Range('E48').Value = '  +4.05%  '
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('E51').Value = '  -0.98%  '
